$d = $word.ActiveDocument

# Locate "Anna Nikolova" and, within it, the "Nikolova" surname, so the
# insertion point for the new middle name is found robustly rather than
# via hard-coded character offsets.
$full = $d.Content
$full.Find.Execute("Anna Nikolova") | Out-Null
$surname = $d.Range($full.Start, $full.End)
$surname.Find.Execute("Nikolova") | Out-Null
$insertAt = $surname.Start

# Insert the new middle name text right before "Nikolova".
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertBefore("Dimitrova ")

# Force Word to keep "Anna ", "Dimitrova " and "Nikolova" as distinct runs
# (otherwise identically-formatted adjacent runs get coalesced back into a
# single run) by toggling a character property on the new text on and back
# off again.
$newRange = $d.Range($insertAt, $insertAt + 10)
$newRange.Font.Bold = $true
$newRange.Font.Bold = $false

# Move the "_GoBack" bookmark so it sits between "Dimitrova " and
# "Nikolova", matching its original position relative to the last edit.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($insertAt + 10, $insertAt + 10))
